$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (Fecha/D, Volumen/J, Precio minimo/K, Precio maximo/L,
# Precio promedio ponderado/M, Precio $/Kg P) after re-shuffling the weekly data.
$data = @{
    2  = @{ D = 44964; J = 300; K = 20000; L = 21000; M = 20500; P = 1139 }
    3  = @{ D = 45068; J = 400; K = 16000; L = 17000; M = 16500; P = 917 }
    4  = @{ D = 44557; J = 400; K = 13000; L = 14000; M = 13500; P = 750 }
    5  = @{ D = 44984; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
    6  = @{ D = 44977; J = 400; K = 16500; L = 17000; M = 16750; P = 931 }
    7  = @{ D = 44547; J = 200; K = 13000; L = 14000; M = 13500; P = 750 }
    8  = @{ D = 44568; J = 500; K = 15000; L = 16000; M = 15500; P = 861 }
    9  = @{ D = 45005; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
    10 = @{ D = 44998; J = 320; K = 17000; L = 18000; M = 17500; P = 972 }
    11 = @{ D = 44960; J = 400; K = 19500; L = 20000; M = 19750; P = 1097 }
    12 = @{ D = 44957; J = 400; K = 21000; L = 22000; M = 21500; P = 1194 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
